$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Argentina, sigmoid 0) - updated values
$ws.Range("C2").Value = 1948
$ws.Range("D2").Value = 76
$ws.Range("E2").Value = 0.1
$ws.Range("F2").Value = 44.2

# Row 3 (sigmoid 1) - updated values
$ws.Range("C3").Value = 3342
$ws.Range("D3").Value = 217
$ws.Range("E3").Value = 0.244
$ws.Range("F3").Value = 33.2

# Row 4 (sigmoid 2) - now new values (previously held the old row4 payload)
$ws.Range("C4").Value = 689
$ws.Range("D4").Value = 324
$ws.Range("E4").Value = 0.1
$ws.Range("F4").Value = 21.5

# Row 5 (new, sigmoid 3) - carries the old row4 C value plus new D/E/F
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 1736
$ws.Range("D5").Value = 418
$ws.Range("E5").Value = 0.1
$ws.Range("F5").Value = 28.4

# Row 6 (new, sigmoid 4)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 2068
$ws.Range("D6").Value = 490
$ws.Range("E6").Value = 0.1
$ws.Range("F6").Value = 40.1

# Match formatting of the other "Sigmoid Number" column cells (bold, bordered,
# centered style) onto the two newly-added rows in column A.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
